# Commit: "addition of static and global variables"
#
# The checklist's "Grade We Think We Will Get" column (D) is filled in for
# the rows the student now credits themselves for, and the running total in
# D33 (=SUM(D3:D32)) recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value  = 1     # Math / Technical Knowledge from other areas
$ws.Range("D10").Value = 1     # Resource Usages
$ws.Range("D11").Value = 1     # Write, compile, and include your own libraries
$ws.Range("D13").Value = 0.5   # Global and Static variables
$ws.Range("D16").Value = 0.5   # Functions OR Recursions
$ws.Range("D18").Value = 1     # Pointers
$ws.Range("D23").Value = 1     # Containers
$ws.Range("D29").Value = 1     # GUI

# Move the saved cursor/selection to match the author's final position.
$ws.Range("E17").Select() | Out-Null
